$wb = $excel.ActiveWorkbook

$sheetCache = @{}
function Get-Sheet($name) {
    if (-not $sheetCache.ContainsKey($name)) {
        $sheetCache[$name] = $wb.Worksheets.Item($name)
    }
    return $sheetCache[$name]
}

# ARM!row17
$ws = Get-Sheet "ARM"
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# ARM!row32
$ws = Get-Sheet "ARM"
$ws.Range("H32").Value = 301
$ws.Range("I32").Value = 301
$ws.Range("K32").Value = 301
$ws.Range("M32").Value = -14

# ARM!row45
$ws = Get-Sheet "ARM"
$ws.Range("H45").Value = 1914.8
$ws.Range("I45").Value = 1608.1666
$ws.Range("K45").Value = 1608.1666
$ws.Range("M45").Value = -1231.1666

# BSM!row5
$ws = Get-Sheet "BSM"
$ws.Range("H5").Value = 259.6
$ws.Range("I5").Value = 308.25
$ws.Range("J5").Value = 65
$ws.Range("K5").Value = 308.25
$ws.Range("L5").Value = 65
$ws.Range("M5").Value = -195.25
$ws.Range("N5").Value = -291

# BSM!row7
$ws = Get-Sheet "BSM"
$ws.Range("H7").Value = 386
$ws.Range("I7").Value = 322
$ws.Range("J7").Value = 450
$ws.Range("K7").Value = 322
$ws.Range("L7").Value = 450
$ws.Range("M7").Value = -209
$ws.Range("N7").Value = -676

# BSM!row37
$ws = Get-Sheet "BSM"
$ws.Range("H37").Value = 775
$ws.Range("I37").Value = 50
$ws.Range("K37").Value = 50
$ws.Range("M37").Value = 87

# BSM!row46
$ws = Get-Sheet "BSM"
$ws.Range("H46").Value = 12000
$ws.Range("I46").Value = 12000
$ws.Range("K46").Value = 12000
$ws.Range("M46").Value = -11702

# BSM!row99
$ws = Get-Sheet "BSM"
$ws.Range("H99").Value = 2124.25
$ws.Range("I99").Value = 1499.5
$ws.Range("J99").Value = 2749
$ws.Range("K99").Value = 1499.5
$ws.Range("L99").Value = 2749
$ws.Range("M99").Value = -1.5
$ws.Range("N99").Value = -5745

# BSM!row107
$ws = Get-Sheet "BSM"
$ws.Range("H107").Value = 4455
$ws.Range("I107").Value = 4455
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4455
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2535
$ws.Range("N107").ClearContents()

# CRP!row7
$ws = Get-Sheet "CRP"
$ws.Range("H7").Value = 93.80851
$ws.Range("I7").Value = 159.06667
$ws.Range("J7").Value = 63.21875
$ws.Range("K7").Value = 159.06667
$ws.Range("L7").Value = 63.21875
$ws.Range("M7").Value = -46.06666999999999
$ws.Range("N7").Value = -289.21875

# CRP!row16
$ws = Get-Sheet "CRP"
$ws.Range("H16").Value = 864.4
$ws.Range("I16").Value = 858.5
$ws.Range("K16").Value = 858.5
$ws.Range("M16").Value = -571.5

# CRP!row17
$ws = Get-Sheet "CRP"
$ws.Range("H17").Value = 900
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# CRP!row31
$ws = Get-Sheet "CRP"
$ws.Range("H31").Value = 1999.5
$ws.Range("I31").Value = 1999.5
$ws.Range("K31").Value = 1999.5
$ws.Range("M31").Value = -1704.5

# CRP!row34
$ws = Get-Sheet "CRP"
$ws.Range("H34").Value = 1999.5
$ws.Range("I34").Value = 1999.5
$ws.Range("K34").Value = 1999.5
$ws.Range("M34").Value = -1797.5

# CRP!row39
$ws = Get-Sheet "CRP"
$ws.Range("H39").Value = 2500
$ws.Range("J39").Value = 2500
$ws.Range("L39").Value = 2500
$ws.Range("N39").Value = -3282

# CRP!row49
$ws = Get-Sheet "CRP"
$ws.Range("H49").Value = 2500
$ws.Range("J49").Value = 2500
$ws.Range("L49").Value = 2500
$ws.Range("N49").Value = -2864

# CRP!row113
$ws = Get-Sheet "CRP"
$ws.Range("H113").Value = 864.4
$ws.Range("I113").Value = 858.5
$ws.Range("K113").Value = 858.5
$ws.Range("M113").Value = 1311.5

# CRP!row122
$ws = Get-Sheet "CRP"
$ws.Range("H122").Value = 1275.4445
$ws.Range("J122").Value = 990
$ws.Range("L122").Value = 2970
$ws.Range("N122").Value = -7870

# CRP!row134
$ws = Get-Sheet "CRP"
$ws.Range("H134").Value = 1999
$ws.Range("I134").Value = 1999
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5997
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3462
$ws.Range("N134").ClearContents()

# CUL!row23
$ws = Get-Sheet "CUL"
$ws.Range("H23").Value = 386.75
$ws.Range("J23").Value = 499
$ws.Range("L23").Value = 1497
$ws.Range("N23").Value = -1967

# CUL!row75
$ws = Get-Sheet "CUL"
$ws.Range("H75").Value = 198
$ws.Range("I75").Value = 198
$ws.Range("K75").Value = 594
$ws.Range("M75").Value = 404

# CUL!row78
$ws = Get-Sheet "CUL"
$ws.Range("H78").Value = 198
$ws.Range("I78").Value = 198
$ws.Range("K78").Value = 1782
$ws.Range("M78").Value = 3210

# CUL!row117
$ws = Get-Sheet "CUL"
$ws.Range("H117").Value = 360
$ws.Range("J117").Value = 386
$ws.Range("L117").Value = 1158
$ws.Range("N117").Value = -8042

# GSM!row44
$ws = Get-Sheet "GSM"
$ws.Range("H44").Value = 25000
$ws.Range("J44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("N44").Value = -26192

# GSM!row47
$ws = Get-Sheet "GSM"
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# GSM!row49
$ws = Get-Sheet "GSM"
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# GSM!row59
$ws = Get-Sheet "GSM"
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -21166

# GSM!row98
$ws = Get-Sheet "GSM"
$ws.Range("H98").Value = 18000
$ws.Range("J98").Value = 18000
$ws.Range("L98").Value = 18000
$ws.Range("N98").Value = -23990

# GSM!row122
$ws = Get-Sheet "GSM"
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# LTW!row9
$ws = Get-Sheet "LTW"
$ws.Range("H9").Value = 821.3333
$ws.Range("I9").Value = 575
$ws.Range("J9").Value = 944.5
$ws.Range("K9").Value = 575
$ws.Range("L9").Value = 944.5
$ws.Range("M9").Value = -351
$ws.Range("N9").Value = -1392.5

# LTW!row30
$ws = Get-Sheet "LTW"
$ws.Range("H30").Value = 765.4
$ws.Range("I30").Value = 765.4
$ws.Range("K30").Value = 765.4
$ws.Range("M30").Value = -657.4

# LTW!row35
$ws = Get-Sheet "LTW"
$ws.Range("H35").Value = 4187.3335
$ws.Range("I35").Value = 4187.3335
$ws.Range("K35").Value = 4187.3335
$ws.Range("M35").Value = -3851.3335

# LTW!row39
$ws = Get-Sheet "LTW"
$ws.Range("H39").Value = 3525
$ws.Range("I39").Value = 3525
$ws.Range("K39").Value = 3525
$ws.Range("M39").Value = -3065

# LTW!row58
$ws = Get-Sheet "LTW"
$ws.Range("H58").Value = 6942
$ws.Range("I58").Value = 6942
$ws.Range("K58").Value = 6942
$ws.Range("M58").Value = -6682

# WVR!row96
$ws = Get-Sheet "WVR"
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# WVR!row136
$ws = Get-Sheet "WVR"
$ws.Range("H136").Value = 1650.3
$ws.Range("I136").Value = 1812.875
$ws.Range("K136").Value = 5438.625
$ws.Range("M136").Value = -2888.625
